$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Range("F4").Value = 3176
$wb.Worksheets.Item(1).Range("F9").Value = 7450
$wb.Worksheets.Item(1).Range("F19").Value = 688
$wb.Worksheets.Item(1).Range("F23").Value = 586
$wb.Worksheets.Item(1).Range("F27").Value = 978
$wb.Worksheets.Item(1).Range("F28").Value = 89
$wb.Worksheets.Item(1).Range("F29").Value = 5809
$wb.Worksheets.Item(1).Range("F31").Value = 4014
$wb.Worksheets.Item(1).Range("F33").Value = 202
$wb.Worksheets.Item(1).Range("F43").Value = 310
$wb.Worksheets.Item(1).Range("F45").Value = 892
$wb.Worksheets.Item(1).Range("F46").Value = 447
$wb.Worksheets.Item(1).Range("F49").Value = 210
$wb.Worksheets.Item(2).Range("F17").Value = 95
$wb.Worksheets.Item(2).Range("F28").Value = 2728
$wb.Worksheets.Item(2).Range("F29").Value = 2728
$wb.Worksheets.Item(3).Range("F6").Value = 1891
$wb.Worksheets.Item(3).Range("F8").Value = 2955
$wb.Worksheets.Item(3).Range("F10").Value = 1178
$wb.Worksheets.Item(3).Range("F14").Value = 8332
$wb.Worksheets.Item(3).Range("F15").Value = 644
$wb.Worksheets.Item(4).Range("F4").Value = 3176
$wb.Worksheets.Item(4).Range("F5").Value = 1891
$wb.Worksheets.Item(4).Range("F7").Value = 2955
$wb.Worksheets.Item(4).Range("F8").Value = 7450
$wb.Worksheets.Item(4).Range("F10").Value = 1178
$wb.Worksheets.Item(4).Range("F16").Value = 688
$wb.Worksheets.Item(4).Range("F24").Value = 586
$wb.Worksheets.Item(4).Range("F29").Value = 978
$wb.Worksheets.Item(4).Range("F30").Value = 89
$wb.Worksheets.Item(4).Range("F31").Value = 5809
$wb.Worksheets.Item(4).Range("F33").Value = 4014
$wb.Worksheets.Item(4).Range("F44").Value = 310
$wb.Worksheets.Item(4).Range("F47").Value = 447
$wb.Worksheets.Item(4).Range("F49").Value = 2728
